# March 24 update 3
# Adds three new columns (N: renewd, O: PlanID, P: iteration) to Sheet1,
# filling the data rows (2-25) with "before", "502-0147678", 14 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("N1").Value = "renewd"
$ws.Range("O1").Value = "PlanID"
$ws.Range("P1").Value = "iteration"

# Match the header style used by the existing header row (B1:M1):
# bold font, thin border on all sides, centered horizontally, top-aligned vertically
$hdr = $ws.Range("N1:P1")
$hdr.Borders.LineStyle = 1
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Fill the new columns for each data row
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 14).Value = "before"
    $ws.Cells.Item($r, 15).Value = "502-0147678"
    $ws.Cells.Item($r, 16).Value = 14
}
